# Apply the edits described by the diff:
# 1. Workbook window height: 7860 -> 5088
# 2. Shared string email text: zeinabragab5new@mailinator.com -> zeinabtest77@mailinator.com
# 3. Sheet selection: C2 -> C1
# 4. Column C width: 29.109375 -> 33.44140625

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Resize the workbook window (best-effort; mirrors the saved bookViews/workbookView size)
$win = $excel.ActiveWindow
$win.Height = 5088
$win.Width = 16860

# 2. Update the email text held in cell C1 (shared string), leaving the hyperlink target untouched
$ws.Range("C1").Value2 = "zeinabtest77@mailinator.com"

# 3. Move/select cell C1 (updates the sheetView's <selection> element)
[void]$ws.Range("C1").Select()

# 4. Widen column C
$ws.Columns.Item(3).ColumnWidth = 33.44140625
